$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target (after) values per row for columns D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# S (Precio $/Kg). These came out of a re-sort/refresh of the weekly
# fruit/vegetable price data; only these six columns move between rows,
# all other columns (A,B,C,E-L,Q,R,T) stay put for a given row.


$ws.Cells.Item(2, 4).Value2 = 44400
$ws.Cells.Item(2, 13).Value2 = 45
$ws.Cells.Item(2, 14).Value2 = 20000
$ws.Cells.Item(2, 15).Value2 = 20000
$ws.Cells.Item(2, 16).Value2 = 20000
$ws.Cells.Item(2, 19).Value2 = 1000

$ws.Cells.Item(3, 4).Value2 = 44305
$ws.Cells.Item(3, 13).Value2 = 20
$ws.Cells.Item(3, 14).Value2 = 22000
$ws.Cells.Item(3, 15).Value2 = 22000
$ws.Cells.Item(3, 16).Value2 = 22000
$ws.Cells.Item(3, 19).Value2 = 1100

$ws.Cells.Item(4, 4).Value2 = 44406
$ws.Cells.Item(4, 13).Value2 = 20
$ws.Cells.Item(4, 14).Value2 = 20000
$ws.Cells.Item(4, 15).Value2 = 20000
$ws.Cells.Item(4, 16).Value2 = 20000
$ws.Cells.Item(4, 19).Value2 = 1000

$ws.Cells.Item(5, 4).Value2 = 44376
$ws.Cells.Item(5, 13).Value2 = 38
$ws.Cells.Item(5, 14).Value2 = 20000
$ws.Cells.Item(5, 15).Value2 = 20000
$ws.Cells.Item(5, 16).Value2 = 20000
$ws.Cells.Item(5, 19).Value2 = 1000

$ws.Cells.Item(6, 4).Value2 = 44382
$ws.Cells.Item(6, 13).Value2 = 24
$ws.Cells.Item(6, 14).Value2 = 20000
$ws.Cells.Item(6, 15).Value2 = 20000
$ws.Cells.Item(6, 16).Value2 = 20000
$ws.Cells.Item(6, 19).Value2 = 1000

$ws.Cells.Item(7, 4).Value2 = 44294
$ws.Cells.Item(7, 13).Value2 = 25
$ws.Cells.Item(7, 14).Value2 = 25000
$ws.Cells.Item(7, 15).Value2 = 25000
$ws.Cells.Item(7, 16).Value2 = 25000
$ws.Cells.Item(7, 19).Value2 = 1250

$ws.Cells.Item(8, 4).Value2 = 44403
$ws.Cells.Item(8, 13).Value2 = 50
$ws.Cells.Item(8, 14).Value2 = 20000
$ws.Cells.Item(8, 15).Value2 = 20000
$ws.Cells.Item(8, 16).Value2 = 20000
$ws.Cells.Item(8, 19).Value2 = 1000

$ws.Cells.Item(9, 4).Value2 = 44413
$ws.Cells.Item(9, 13).Value2 = 45
$ws.Cells.Item(9, 14).Value2 = 20000
$ws.Cells.Item(9, 15).Value2 = 20000
$ws.Cells.Item(9, 16).Value2 = 20000
$ws.Cells.Item(9, 19).Value2 = 1000

$ws.Cells.Item(10, 4).Value2 = 44291
$ws.Cells.Item(10, 13).Value2 = 70
$ws.Cells.Item(10, 14).Value2 = 25000
$ws.Cells.Item(10, 15).Value2 = 25000
$ws.Cells.Item(10, 16).Value2 = 25000
$ws.Cells.Item(10, 19).Value2 = 1250

$ws.Cells.Item(11, 4).Value2 = 44292
$ws.Cells.Item(11, 13).Value2 = 30
$ws.Cells.Item(11, 14).Value2 = 25000
$ws.Cells.Item(11, 15).Value2 = 25000
$ws.Cells.Item(11, 16).Value2 = 25000
$ws.Cells.Item(11, 19).Value2 = 1250

$ws.Cells.Item(12, 4).Value2 = 44445
$ws.Cells.Item(12, 13).Value2 = 45
$ws.Cells.Item(12, 14).Value2 = 20000
$ws.Cells.Item(12, 15).Value2 = 20000
$ws.Cells.Item(12, 16).Value2 = 20000
$ws.Cells.Item(12, 19).Value2 = 1000

$ws.Cells.Item(13, 4).Value2 = 44300
$ws.Cells.Item(13, 13).Value2 = 45
$ws.Cells.Item(13, 14).Value2 = 22000
$ws.Cells.Item(13, 15).Value2 = 22000
$ws.Cells.Item(13, 16).Value2 = 22000
$ws.Cells.Item(13, 19).Value2 = 1100

$ws.Cells.Item(14, 4).Value2 = 44448
$ws.Cells.Item(14, 13).Value2 = 30
$ws.Cells.Item(14, 14).Value2 = 22000
$ws.Cells.Item(14, 15).Value2 = 22000
$ws.Cells.Item(14, 16).Value2 = 22000
$ws.Cells.Item(14, 19).Value2 = 1100

$ws.Cells.Item(15, 4).Value2 = 44385
$ws.Cells.Item(15, 13).Value2 = 36
$ws.Cells.Item(15, 14).Value2 = 20000
$ws.Cells.Item(15, 15).Value2 = 20000
$ws.Cells.Item(15, 16).Value2 = 20000
$ws.Cells.Item(15, 19).Value2 = 1000

$ws.Cells.Item(16, 4).Value2 = 44377
$ws.Cells.Item(16, 13).Value2 = 25
$ws.Cells.Item(16, 14).Value2 = 20000
$ws.Cells.Item(16, 15).Value2 = 20000
$ws.Cells.Item(16, 16).Value2 = 20000
$ws.Cells.Item(16, 19).Value2 = 1000

$ws.Cells.Item(17, 4).Value2 = 44307
$ws.Cells.Item(17, 13).Value2 = 30
$ws.Cells.Item(17, 14).Value2 = 22000
$ws.Cells.Item(17, 15).Value2 = 22000
$ws.Cells.Item(17, 16).Value2 = 22000
$ws.Cells.Item(17, 19).Value2 = 1100

$ws.Cells.Item(18, 4).Value2 = 44301
$ws.Cells.Item(18, 13).Value2 = 38
$ws.Cells.Item(18, 14).Value2 = 22000
$ws.Cells.Item(18, 15).Value2 = 22000
$ws.Cells.Item(18, 16).Value2 = 22000
$ws.Cells.Item(18, 19).Value2 = 1100

$ws.Cells.Item(19, 4).Value2 = 44298
$ws.Cells.Item(19, 13).Value2 = 65
$ws.Cells.Item(19, 14).Value2 = 22000
$ws.Cells.Item(19, 15).Value2 = 22000
$ws.Cells.Item(19, 16).Value2 = 22000
$ws.Cells.Item(19, 19).Value2 = 1100

$ws.Cells.Item(20, 4).Value2 = 44389
$ws.Cells.Item(20, 13).Value2 = 20
$ws.Cells.Item(20, 14).Value2 = 20000
$ws.Cells.Item(20, 15).Value2 = 20000
$ws.Cells.Item(20, 16).Value2 = 20000
$ws.Cells.Item(20, 19).Value2 = 1000
